$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-7
# from serial date 45208 (2023-10-09) to 45212 (2023-10-13)
$ws.Range("C2:C7").Value = 45212
